# Update the division-problem answers in the single table of the document.
# Cells are addressed by (row, column) rather than by text search, because
# some old values are duplicated (e.g. "48÷3=16, 0" appears twice) and must
# be replaced with different new values depending on position.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map: Word table row (1-based) -> array of new column values (columns 1..5)
$updates = @{
    1  = @("54÷6=9, 0", "26÷2=13, 0", "21÷9=2, 3", "31÷5=6, 1", "76÷7=10, 6")
    5  = @("28÷5=5, 3", "53÷9=5, 8", "52÷5=10, 2", "58÷3=19, 1", "83÷5=16, 3")
    9  = @("13÷7=1, 6", "67÷2=33, 1", "62÷8=7, 6", "16÷2=8, 0", "18÷8=2, 2")
    13 = @("79÷5=15, 4", "59÷4=14, 3", "32÷4=8, 0", "86÷3=28, 2", "57÷5=11, 2")
    17 = @("62÷8=7, 6", "39÷9=4, 3", "35÷9=3, 8", "22÷4=5, 2", "98÷9=10, 8")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
